$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 15.63
$ws.Range("E2").Value = 64.09999999999999
$ws.Range("F2").Value = 4.46
$ws.Range("N2").Value = 53.71147335634279

# Row 3
$ws.Range("D3").Value = 92425.28999999999
$ws.Range("E3").Value = 62.2
$ws.Range("F3").Value = 1.73
$ws.Range("N3").Value = 53.71147335634279

# Row 4
$ws.Range("D4").Value = 275.92
$ws.Range("E4").Value = 46.9
$ws.Range("F4").Value = 4.13
$ws.Range("N4").Value = 53.71147335634279

# Row 5
$ws.Range("D5").Value = 12.46
$ws.Range("E5").Value = 47.6
$ws.Range("F5").Value = 12.14
$ws.Range("N5").Value = 53.71147335634279

# Row 6
$ws.Range("D6").Value = 187.1
$ws.Range("E6").Value = 39.8
$ws.Range("F6").Value = 6.53
$ws.Range("G6").Value = 30
$ws.Range("K6").Value = 35.5
$ws.Range("N6").Value = 53.71147335634279
